$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from an existing header cell
# into the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-40).
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 65
    $ws.Cells.Item($r, 32).Value = 0
}
